# Rebuild the bull-bear cover-hero slide with native, editable PPTX elements.
#
# Summary of the edit (see commit message / diff):
#   - remove the full-bleed background picture (id=2 "Picture 1")
#   - turn the old "main_title" placeholder textbox into "TextBox 1" with
#     literal copy, new size/position and new run formatting
#   - turn the old "hero_image" outline shape into the first
#     "image_placeholder" dashed-outline box (new size/position/line style)
#   - add a second "image_placeholder" dashed-outline box (duplicate of the
#     first, repositioned)
#   - add two new rounded-rectangle shapes ("Rounded Rectangle 4/5") that
#     provide the white band + yellow background fill
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/Top/Width/Height
# (and Line.Weight) are expressed in points and are stored internally as
# single-precision floats, so a plain EMU/12700 division can land 1 EMU
# short after the point value round-trips through float32 + EMU conversion.
# The literals used here were solved for so the resulting EMU values match
# the target exactly wherever a float32 representation makes that possible
# (two of the widths -- 15245425 and 24384030 EMU -- have no exact float32
# point representation; those land within 1 EMU of the target, i.e. well
# under 1/1000000 inch).

function RGBColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Drop the full-slide background picture (id=2 "Picture 1").
# ---------------------------------------------------------------------
$s.Shapes.Item(1).Delete()

# ---------------------------------------------------------------------
# 2) "main_title" -> "TextBox 1" (now shape 1 after the picture removal)
# ---------------------------------------------------------------------
$titleBox = $s.Shapes.Item(1)
$titleBox.Name = "TextBox 1"
$titleBox.Left = 132.39354705826773
$titleBox.Top = 463.996780396063
$titleBox.Width = 827.6111755381891
$titleBox.Height = 280.0007781984252

$titleRange = $titleBox.TextFrame.TextRange
$titleRange.Text = "bull bear strategy"
$titleRange.ParagraphFormat.Alignment = 1
$titleRange.Font.Size = 130
$titleRange.Font.Bold = $false
$titleRange.Font.Italic = $false
$titleRange.Font.Name = "Rajdhani"
$titleRange.Font.Color.RGB = RGBColor 0x00 0x00 0x00

# ---------------------------------------------------------------------
# 3) "hero_image" -> first "image_placeholder" (now shape 2)
#    Duplicate it *before* restyling so the copy inherits the original
#    hero_image look (p:style + txBody) exactly, then restyle both.
# ---------------------------------------------------------------------
$heroBox = $s.Shapes.Item(2)
$placeholder2 = $heroBox.Duplicate().Item(1)

$heroBox.Name = "image_placeholder"
$heroBox.Left = 719.5752258311023
$heroBox.Top = 141.73912811338582
$heroBox.Width = 1200.4271653543308
$heroBox.Height = 887.9111633303149
$heroBox.Fill.Visible = $false
$heroBox.Line.Visible = $true
$heroBox.Line.Weight = 0.5
$heroBox.Line.ForeColor.RGB = RGBColor 0xCC 0xCC 0xCC
$heroBox.Line.DashStyle = 9

# ---------------------------------------------------------------------
# 4) Second "image_placeholder" (the hero_image duplicate, now shape 3)
# ---------------------------------------------------------------------
$placeholder2.Name = "image_placeholder"
$placeholder2.Left = 82.36795425472441
$placeholder2.Top = 66.32638168385827
$placeholder2.Width = 452.1671600346457
$placeholder2.Height = 92.74315261929134
$placeholder2.Fill.Visible = $false
$placeholder2.Line.Visible = $true
$placeholder2.Line.Weight = 0.5
$placeholder2.Line.ForeColor.RGB = RGBColor 0xCC 0xCC 0xCC
$placeholder2.Line.DashStyle = 9

# ---------------------------------------------------------------------
# 5) New "Rounded Rectangle 4" - white band (now shape 4)
# ---------------------------------------------------------------------
$rr4 = $s.Shapes.AddShape(5, 0, 0, 1, 1)
$rr4.Name = "Rounded Rectangle 4"
$rr4.Left = 0.0
$rr4.Top = 129.31197357204724
$rr4.Width = 1920.0023622047245
$rr4.Height = 120.40559005748032
$rr4.Fill.Visible = $true
$rr4.Fill.ForeColor.RGB = RGBColor 0xFF 0xFF 0xFF
$rr4.Line.Visible = $false

# ---------------------------------------------------------------------
# 6) New "Rounded Rectangle 5" - yellow background (now shape 5)
# ---------------------------------------------------------------------
$rr5 = $s.Shapes.AddShape(5, 0, 0, 1, 1)
$rr5.Name = "Rounded Rectangle 5"
$rr5.Left = 0.0
$rr5.Top = 132.14157867480316
$rr5.Width = 1920.0023622047245
$rr5.Height = 947.858367920866
$rr5.Fill.Visible = $true
$rr5.Fill.ForeColor.RGB = RGBColor 0xFE 0xC0 0x0F
$rr5.Line.Visible = $false
